# s1cDNASample_hbrown_07.22.19.xlsx - "further cleaning to metadata"
#
# 1. Column H (kit/reagent lot number) sample value is corrected from
#    "E7760" to "E7420" for every data row (H2:H27).
# 2. The font metadata backing that column's style is cleaned up to the
#    plain Arial 10 / generic-family font used elsewhere in the sheet
#    (this is what produces the new, distinct cell style for H2:H27).
# 3. The active selection is moved from E2:E27 onto the corrected column,
#    H2:H27.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("H2:H27")

# Fix the sample value.
$rng.Value = "E7420"

# Clean up the font metadata for this range.
$rng.Font.Name = "Arial"
$rng.Font.Size = 10
$rng.Font.Family = 0

# Move the selection onto the corrected column.
$rng.Select()
